$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "58.684.86"
$ws.Cells.Item(2, 5).Value = "  -2.85%  "

$ws.Cells.Item(3, 4).Value = "2.724.87"
$ws.Cells.Item(3, 5).Value = "  -5.75%  "

$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 5).Value = "  +0.15%  "

$ws.Cells.Item(5, 4).Value = "503.78"
$ws.Cells.Item(5, 5).Value = "  -4.43%  "

$ws.Cells.Item(6, 4).Value = "141.05"
$ws.Cells.Item(6, 5).Value = "  -0.78%  "

$ws.Cells.Item(7, 4).Value = "0.998"
$ws.Cells.Item(7, 5).Value = "  -0.13%  "

$ws.Cells.Item(8, 4).Value = "'0.530"
$ws.Cells.Item(8, 5).Value = "  -3.85%  "

$ws.Cells.Item(9, 4).Value = "2.739.40"
$ws.Cells.Item(9, 5).Value = "  -5.26%  "

$ws.Cells.Item(10, 4).Value = "6.06"
$ws.Cells.Item(10, 5).Value = "  +2.64%  "

$ws.Cells.Item(11, 5).Value = "  -2.30%  "

$ws.Cells.Item(12, 4).Value = "0.347"
$ws.Cells.Item(12, 5).Value = "  -2.89%  "

$ws.Cells.Item(13, 4).Value = "0.126"
$ws.Cells.Item(13, 5).Value = "  +1.12%  "

$ws.Cells.Item(14, 4).Value = "3.203.65"
$ws.Cells.Item(14, 5).Value = "  -5.41%  "

$ws.Cells.Item(15, 4).Value = "58.834.88"
$ws.Cells.Item(15, 5).Value = "  -2.58%  "

$ws.Cells.Item(16, 4).Value = "21.68"
$ws.Cells.Item(16, 5).Value = "  -3.77%  "

$ws.Cells.Item(17, 4).Value = "2.740.13"
$ws.Cells.Item(17, 5).Value = "  -4.67%  "

$ws.Cells.Item(18, 4).Value = "0.0000135"
$ws.Cells.Item(18, 5).Value = "  -3.88%  "

$ws.Cells.Item(19, 4).Value = "4.76"
$ws.Cells.Item(19, 5).Value = "  -3.37%  "

$ws.Cells.Item(20, 4).Value = "10.98"
$ws.Cells.Item(20, 5).Value = "  -4.93%  "

$ws.Cells.Item(21, 4).Value = "'343.50"
$ws.Cells.Item(21, 5).Value = "  -4.86%  "

$ws.Cells.Item(22, 4).Value = "6.25"
$ws.Cells.Item(22, 5).Value = "  -4.14%  "

$ws.Cells.Item(23, 5).Value = "  -0.55%  "

$ws.Cells.Item(24, 5).Value = "  -0.46%  "

$ws.Cells.Item(25, 4).Value = "63.01"
$ws.Cells.Item(25, 5).Value = "  -0.42%  "

$ws.Cells.Item(26, 4).Value = "0.426"
$ws.Cells.Item(26, 5).Value = "  -5.13%  "

$ws.Cells.Item(27, 5).Value = "  -4.99%  "

$ws.Cells.Item(28, 4).Value = "0.995"
$ws.Cells.Item(28, 5).Value = "  -0.42%  "

$ws.Cells.Item(29, 4).Value = "7.52"
$ws.Cells.Item(29, 5).Value = "  -3.33%  "

$ws.Cells.Item(30, 4).Value = "0.0₃0830"
$ws.Cells.Item(30, 5).Value = "  -2.83%  "

$ws.Cells.Item(31, 5).Value = "  -0.12%  "

$ws.Cells.Item(32, 2).Value = "EthereumClassic"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(32, 4).Value = "19.16"
$ws.Cells.Item(32, 5).Value = "  -1.49%  "

$ws.Cells.Item(33, 2).Value = "PancakeSwap"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(33, 4).Value = "'1.60"
$ws.Cells.Item(33, 5).Value = "  -4.07%  "

$ws.Cells.Item(34, 4).Value = "151.74"
$ws.Cells.Item(34, 5).Value = "  +1.40%  "

$ws.Cells.Item(35, 4).Value = "5.43"
$ws.Cells.Item(35, 5).Value = "  -2.19%  "

$ws.Cells.Item(36, 5).Value = "  -3.15%  "

$ws.Cells.Item(37, 4).Value = "0.949"
$ws.Cells.Item(37, 5).Value = "  -4.02%  "

$ws.Cells.Item(38, 5).Value = "  -5.76%  "

$ws.Cells.Item(39, 4).Value = "35.91"
$ws.Cells.Item(39, 5).Value = "  -4.89%  "

$ws.Cells.Item(40, 4).Value = "1.39"
$ws.Cells.Item(40, 5).Value = "  -6.62%  "

$ws.Cells.Item(41, 4).Value = "3.55"
$ws.Cells.Item(41, 5).Value = "  -2.82%  "

$ws.Cells.Item(42, 4).Value = "2.189.37"
$ws.Cells.Item(42, 5).Value = "  -5.70%  "

$ws.Cells.Item(43, 4).Value = "0.0559"
$ws.Cells.Item(43, 5).Value = "  -2.13%  "

$ws.Cells.Item(44, 4).Value = "0.998"
$ws.Cells.Item(44, 5).Value = "  +0.12%  "

$ws.Cells.Item(45, 4).Value = "0.604"
$ws.Cells.Item(45, 5).Value = "  -5.51%  "

$ws.Cells.Item(46, 4).Value = "19.02"
$ws.Cells.Item(46, 5).Value = "  -7.95%  "

$ws.Cells.Item(47, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Cells.Item(47, 4).Value = "10.37"
$ws.Cells.Item(47, 5).Value = "  +0.31%  "

$ws.Cells.Item(48, 2).Value = "RenderToken"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(48, 4).Value = "4.76"
$ws.Cells.Item(48, 5).Value = "  -5.81%  "

$ws.Cells.Item(49, 5).Value = "  -3.06%  "

$ws.Cells.Item(50, 4).Value = "0.0887"
$ws.Cells.Item(50, 5).Value = "  -4.41%  "

$ws.Cells.Item(51, 4).Value = "18.04"
$ws.Cells.Item(51, 5).Value = "  -1.39%  "
